$wb = $excel.ActiveWorkbook

# --- Sheet "ads" (sheet1.xml): add length-count helper formulas in L/M for
#     each of the 4 ad rows (3, 6, 9, 12). ---
$wsAds = $wb.Worksheets.Item("ads")

$wsAds.Range("L3").Formula = '=LEN(TRIM(D3))-LEN(SUBSTITUTE(D3," ",""))+1'
$wsAds.Range("M3").Formula = '=LEN(TRIM(F3))-LEN(SUBSTITUTE(F3," ",""))+1'

$wsAds.Range("L6").Formula = '=LEN(TRIM(D6))-LEN(SUBSTITUTE(D6," ",""))+1'
$wsAds.Range("M6").Formula = '=LEN(TRIM(F6))-LEN(SUBSTITUTE(F6," ",""))+1'

$wsAds.Range("L9").Formula = '=LEN(TRIM(D9))-LEN(SUBSTITUTE(D9," ",""))+1'
$wsAds.Range("M9").Formula = '=LEN(TRIM(F9))-LEN(SUBSTITUTE(F9," ",""))+1'

$wsAds.Range("L12").Formula = '=LEN(TRIM(D12))-LEN(SUBSTITUTE(D12," ",""))+1'
$wsAds.Range("M12").Formula = '=LEN(TRIM(F12))-LEN(SUBSTITUTE(F12," ",""))+1'

# Selection left on this sheet before moving on (matches the final saved
# cursor position for this tab).
[void]$wsAds.Range("L1:M1048576").Select()

# --- Sheet "blended_personalities" (sheet2.xml): same helper formulas. ---
$wsBlend = $wb.Worksheets.Item("blended_personalities")
[void]$wsBlend.Select()

$wsBlend.Range("L3").Formula = '=LEN(TRIM(D3))-LEN(SUBSTITUTE(D3," ",""))+1'
$wsBlend.Range("M3").Formula = '=LEN(TRIM(F3))-LEN(SUBSTITUTE(F3," ",""))+1'

$wsBlend.Range("L6").Formula = '=LEN(TRIM(D6))-LEN(SUBSTITUTE(D6," ",""))+1'
$wsBlend.Range("M6").Formula = '=LEN(TRIM(F6))-LEN(SUBSTITUTE(F6," ",""))+1'

$wsBlend.Range("L9").Formula = '=LEN(TRIM(D9))-LEN(SUBSTITUTE(D9," ",""))+1'
$wsBlend.Range("M9").Formula = '=LEN(TRIM(F9))-LEN(SUBSTITUTE(F9," ",""))+1'

$wsBlend.Range("L12").Formula = '=LEN(TRIM(D12))-LEN(SUBSTITUTE(D12," ",""))+1'
$wsBlend.Range("M12").Formula = '=LEN(TRIM(F12))-LEN(SUBSTITUTE(F12," ",""))+1'

# This sheet ends up the active tab, scrolled down a bit with L3:M12 selected.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
[void]$wsBlend.Range("L3:M12").Select()
